$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update title/timestamp cell (A1) - refreshed data pull time
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 15 de Junio de 2020 a las 16:46"

# Update country rows: col A = Pais, B = Casos totales, C = Nuevos casos,
# D = Casos activos, E = Recuperados, F = Casos criticos, G = Muertes hoy, H = Muertes
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 2165872
$ws.Cells.Item(4, 3).Value = 3644
$ws.Cells.Item(4, 4).Value = 870077
$ws.Cells.Item(4, 5).Value = 1177912
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 25
$ws.Cells.Item(4, 8).Value = 117883

$ws.Cells.Item(35, 1).Value = "Portugal"
$ws.Cells.Item(35, 2).Value = 37036
$ws.Cells.Item(35, 3).Value = 346
$ws.Cells.Item(35, 4).Value = 22852
$ws.Cells.Item(35, 5).Value = 12664
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 3
$ws.Cells.Item(35, 8).Value = 1520

$ws.Cells.Item(61, 1).Value = "Moldavia"
$ws.Cells.Item(61, 2).Value = 11740
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(61, 4).Value = 6794
$ws.Cells.Item(61, 5).Value = 4538
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 2
$ws.Cells.Item(61, 8).Value = 408

$ws.Cells.Item(77, 1).Value = "Tayikistan"
$ws.Cells.Item(77, 2).Value = 5097
$ws.Cells.Item(77, 3).Value = 62
$ws.Cells.Item(77, 4).Value = 3503
$ws.Cells.Item(77, 5).Value = 1544
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 50

$ws.Cells.Item(78, 1).Value = "Costa de Marfil"
$ws.Cells.Item(78, 2).Value = 5084
$ws.Cells.Item(78, 3).Value = 0
$ws.Cells.Item(78, 4).Value = 2505
$ws.Cells.Item(78, 5).Value = 2534
$ws.Cells.Item(78, 6).Value = 0
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 45

$ws.Cells.Item(131, 1).Value = "Principado de Andorra"
$ws.Cells.Item(131, 2).Value = 853
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(131, 4).Value = 789
$ws.Cells.Item(131, 5).Value = 13
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 51

$ws.Cells.Item(134, 1).Value = "Cabo Verde"
$ws.Cells.Item(134, 2).Value = 759
$ws.Cells.Item(134, 3).Value = 9
$ws.Cells.Item(134, 4).Value = 301
$ws.Cells.Item(134, 5).Value = 452
$ws.Cells.Item(134, 6).Value = 0
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 6

$ws.Cells.Item(148, 1).Value = "Estado de Palestina"
$ws.Cells.Item(148, 2).Value = 501
$ws.Cells.Item(148, 3).Value = 9
$ws.Cells.Item(148, 4).Value = 415
$ws.Cells.Item(148, 5).Value = 83
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(148, 8).Value = 3

$ws.Cells.Item(149, 1).Value = "Liberia"
$ws.Cells.Item(149, 2).Value = 498
$ws.Cells.Item(149, 3).Value = 40
$ws.Cells.Item(149, 4).Value = 221
$ws.Cells.Item(149, 5).Value = 244
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 7).Value = 1
$ws.Cells.Item(149, 8).Value = 33

$ws.Cells.Item(150, 1).Value = "Reunion"
$ws.Cells.Item(150, 2).Value = 495
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 460
$ws.Cells.Item(150, 5).Value = 34
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 1

$ws.Cells.Item(151, 1).Value = "Suazilandia"
$ws.Cells.Item(151, 2).Value = 490
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(151, 4).Value = 249
$ws.Cells.Item(151, 5).Value = 237
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 4

$ws.Cells.Item(152, 1).Value = "Benin"
$ws.Cells.Item(152, 2).Value = 483
$ws.Cells.Item(152, 3).Value = 41
$ws.Cells.Item(152, 4).Value = 232
$ws.Cells.Item(152, 5).Value = 242
$ws.Cells.Item(152, 6).Value = 0
$ws.Cells.Item(152, 7).Value = 3
$ws.Cells.Item(152, 8).Value = 9

$ws.Cells.Item(160, 1).Value = "Birmania"
$ws.Cells.Item(160, 2).Value = 262
$ws.Cells.Item(160, 3).Value = 1
$ws.Cells.Item(160, 4).Value = 175
$ws.Cells.Item(160, 5).Value = 81
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 0
$ws.Cells.Item(160, 8).Value = 6

$ws.Cells.Item(206, 1).Value = "Islas Malvinas"
$ws.Cells.Item(206, 2).Value = 13
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 13
$ws.Cells.Item(206, 5).Value = 0
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0

$ws.Cells.Item(207, 1).Value = "Groenlandia"
$ws.Cells.Item(207, 2).Value = 13
$ws.Cells.Item(207, 3).Value = 0
$ws.Cells.Item(207, 4).Value = 13
$ws.Cells.Item(207, 5).Value = 0
$ws.Cells.Item(207, 6).Value = 0
$ws.Cells.Item(207, 7).Value = 0
$ws.Cells.Item(207, 8).Value = 0

$ws.Cells.Item(208, 1).Value = "Santa Sede"
$ws.Cells.Item(208, 2).Value = 12
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 12
$ws.Cells.Item(208, 5).Value = 0
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 0

$ws.Cells.Item(209, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(209, 2).Value = 12
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 11
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 1

